$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: copy the header formatting from G1 ("sum") so the
# new "Save" column looks like the rest of the header row, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2 for the "Save" column.
$ws.Range("H2").Value = 0
